# Add 2022-Q3 data
#
# 1) Duplicate the existing "2022-Q2" sheet (so we inherit its styles/layout),
#    place the copy right before it, rename to "2022-Q3" and fill it with the
#    new quarter's fund-holding data.
# 2) Insert a new row at the top of the "总计" (summary) sheet for the new
#    2022-Q3 quarter, shifting the existing rows down, and renumber the
#    running index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet from a copy of "2022-Q2"
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2Index = $q2.Index
$q2.Copy($q2, $null)

# The copy is inserted immediately before the source sheet, so it now sits
# at the source sheet's former position.
$q3 = $wb.Worksheets.Item($q2Index)
$q3.Name = "2022-Q3"

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text even when the value looks numeric,
    # then strip the temporary number-format styling back off again.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $q3.Range("B2") "159617"
$q3.Range("C2").Value = "华夏中证智选500价值稳健策略ETF"
Set-TextValue $q3.Range("D2") "2.93"
Set-TextValue $q3.Range("E2") "97.05"
Set-TextValue $q3.Range("F2") "1.44"
Set-TextValue $q3.Range("G2") "0.0422"
$q3.Range("H2").Value = 4

# ---------------------------------------------------------------------------
# Step 2: add the new summary row to the "总计" sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy the number formatting from the row below so column A keeps the same
# style as the other index cells.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$total.Range("B2:D2").ClearFormats()
$total.Application.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.04

# Renumber the running index in column A for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally-active worksheet (unaffected by this change).
$wb.Worksheets.Item("2020-Q4").Activate()
